# Update NATMI LR-pair output (Col18a1-Gpc1) with new TPM-derived values.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New per-cluster "ligand" (sending cluster) inputs -------------------
# Ligand-expressing cells (E) and ligand average expression value (G).
# Ligand detection rate (F) = E / 3 ; Ligand total expression value (H) = G * 3.
$ligandClusters = @("ECs", "FAPs", "MuSCs", "Resolving-Mac")
$ligandE = @{ "ECs" = 3; "FAPs" = 3; "MuSCs" = 3; "Resolving-Mac" = 3 }
$ligandG = @{
    "ECs"           = 40.34291466666667
    "FAPs"          = 41.443863
    "MuSCs"         = 36.642055
    "Resolving-Mac" = 0.07370033333333333
}

# --- New per-cluster "receptor" (target cluster) inputs -------------------
# Receptor average expression value (M); Receptor total expression value (N) = M * 3.
$receptorClusters = @("ECs", "FAPs", "MuSCs", "Resolving-Mac")
$receptorM = @{
    "ECs"           = 0.484733
    "FAPs"          = 7.020353
    "MuSCs"         = 53.289524
    "Resolving-Mac" = 0.4014323333333333
}

# Totals used for the "derived specificity" (share-of-total) columns.
$gTotal = 0
foreach ($c in $ligandClusters) { $gTotal += $ligandG[$c] }
$mTotal = 0
foreach ($c in $receptorClusters) { $mTotal += $receptorM[$c] }

$lastRow = 17
for ($row = 2; $row -le $lastRow; $row++) {
    $sendCluster = $ws.Cells.Item($row, 1).Value2
    $targetCluster = $ws.Cells.Item($row, 4).Value2

    $E = $ligandE[$sendCluster]
    $F = $E / 3
    $G = $ligandG[$sendCluster]
    $H = $G * 3
    $I = $G / $gTotal
    $J = $I

    $M = $receptorM[$targetCluster]
    $N = $M * 3
    $O = $M / $mTotal
    $P = $O

    $Q = $G * $M
    $R = $H * $N
    $S = $I * $O
    $T = $J * $P

    $ws.Cells.Item($row, 5).Value = $E    # E: Ligand-expressing cells
    $ws.Cells.Item($row, 6).Value = $F    # F: Ligand detection rate
    $ws.Cells.Item($row, 7).Value = $G    # G: Ligand average expression value
    $ws.Cells.Item($row, 8).Value = $H    # H: Ligand total expression value
    $ws.Cells.Item($row, 9).Value = $I    # I: Ligand derived specificity of average expression value
    $ws.Cells.Item($row, 10).Value = $J   # J: Ligand derived specificity of total expression value

    $ws.Cells.Item($row, 13).Value = $M   # M: Receptor average expression value
    $ws.Cells.Item($row, 14).Value = $N   # N: Receptor total expression value
    $ws.Cells.Item($row, 15).Value = $O   # O: Receptor derived specificity of average expression value
    $ws.Cells.Item($row, 16).Value = $P   # P: Receptor derived specificity of total expression value

    $ws.Cells.Item($row, 17).Value = $Q   # Q: Edge average expression weight
    $ws.Cells.Item($row, 18).Value = $R   # R: Edge total expression weight
    $ws.Cells.Item($row, 19).Value = $S   # S: Edge average expression derived specificity
    $ws.Cells.Item($row, 20).Value = $T   # T: Edge total expression derived specificity
}
